# Refresh the Price (D) and Volume(1h) (E) columns on the cryptos
# worksheet, mirroring the GitHub Actions scheduled data pull.
#
# Most Price cells hold plain decimal-looking text (e.g. "1.033"); Excel's
# normal type inference would silently convert those to numbers when
# assigned through .Value. To keep them as text (matching the source
# data's inline-string cells) we briefly force a Text number format
# before the assignment, then restore the cell's original ("Normal")
# style so no residual formatting change is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.735.22'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.849.94'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.033'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.029'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4388'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3797'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07410'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8831'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = '1.860.56'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.505'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.686'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07167'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '85.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.037'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009068'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = '27.768.20'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.274'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").Value = '2.088.39'
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.045'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.52'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.994'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.52%  '
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09066'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7744'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.213'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.005'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.571'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.031'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.149'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01973'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05271'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.847'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5177'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1671'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.844'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.717'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("E46").Value = '  +1.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.032'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06576'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.708'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4696'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.93%  '
